$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -22.146
$ws.Range("C4").Value = -12.927
$ws.Range("C5").Value = -12.927
$ws.Range("A6").Value = -21.108
$ws.Range("A7").Value = -21.038
$ws.Range("C8").Value = -12.8
$ws.Range("A16").Value = -21.142
$ws.Range("C16").Value = -12.339
$ws.Range("A20").Value = -21.936
$ws.Range("C22").Value = -12.78
